# updated api key setter
# Applies corrected export values (column B) for a set of rows in Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    33  = 3316022133985
    35  = 34222278.682
    37  = 21412416606.361
    52  = 585395916302.973
    89  = 2788495475.762
    94  = 3266988807.77
    102 = 5704483311.693
    104 = 4417050677.857
    110 = 737116422.265
    131 = 633870713.4859999
    138 = 6036188264.211
    159 = 9535273938
    162 = 11013140616.688
}

foreach ($row in $updates.Keys) {
    $ws.Range("B$row").Value = $updates[$row]
}
